$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# --- Update column D ("eje vertical") values per corrected node coordinates ---
$ws.Cells.Item(153, 4).Value = 15
$ws.Cells.Item(154, 4).Value = 15
$ws.Cells.Item(156, 4).Value = 14
$ws.Cells.Item(157, 4).Value = 14
$ws.Cells.Item(158, 4).Value = 14
$ws.Cells.Item(161, 4).Value = 13
$ws.Cells.Item(162, 4).Value = 13
$ws.Cells.Item(163, 4).Value = 13
$ws.Cells.Item(166, 4).Value = 12
$ws.Cells.Item(167, 4).Value = 12
$ws.Cells.Item(170, 4).Value = 11
$ws.Cells.Item(171, 4).Value = 12
$ws.Cells.Item(172, 4).Value = 11
$ws.Cells.Item(173, 4).Value = 11
$ws.Cells.Item(176, 4).Value = 10
$ws.Cells.Item(177, 4).Value = 10
$ws.Cells.Item(178, 4).Value = 10
$ws.Cells.Item(181, 4).Value = 9
$ws.Cells.Item(182, 4).Value = 9
$ws.Cells.Item(183, 4).Value = 9
$ws.Cells.Item(186, 4).Value = 8
$ws.Cells.Item(187, 4).Value = 8
$ws.Cells.Item(188, 4).Value = 8
$ws.Cells.Item(190, 4).Value = 7
$ws.Cells.Item(191, 4).Value = 7
$ws.Cells.Item(192, 4).Value = 7
$ws.Cells.Item(194, 4).Value = 6
$ws.Cells.Item(195, 4).Value = 6
$ws.Cells.Item(196, 4).Value = 6
$ws.Cells.Item(198, 4).Value = 5
$ws.Cells.Item(199, 4).Value = 5
$ws.Cells.Item(200, 4).Value = 5
$ws.Cells.Item(202, 4).Value = 4
$ws.Cells.Item(203, 4).Value = 4
$ws.Cells.Item(204, 4).Value = 4
$ws.Cells.Item(206, 4).Value = 3
$ws.Cells.Item(207, 4).Value = 3
$ws.Cells.Item(208, 4).Value = 3
$ws.Cells.Item(210, 4).Value = 2
$ws.Cells.Item(211, 4).Value = 2

# --- View changes: freeze top row, set active cell ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G168").Select() | Out-Null
